$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("saldo")

# New row of spec data for agent (BLNC020 / Saldo Deposit Agent)
$ws.Range("B7").Value = "BLNC020"
$ws.Range("A7").Value = "Saldo Deposit Agent"

# Make "saldo" the active sheet/tab and select A8
$ws.Activate()
$ws.Range("A8").Select()
